$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-11, columns B:G (ME, MAE, MSE, RMSE, SE, N)
$data = @(
    @(-0.1365980374615512, 0.5567307316750884, 0.6103753252590118, 0.7812652080177459, 0.7768852086692588, 51),
    @(-0.1983056728170303, 0.5132083419313908, 0.4400880119005959, 0.6633912962201086, 0.6394855011274496, 50),
    @(-0.08011458471284517, 0.50764697151692, 0.518069484399113, 0.7197704386810513, 0.7227105481572865, 49),
    @(-0.1481904045569161, 0.5126676976066116, 0.4729556061377808, 0.6877176790935222, 0.6786684410696205, 48),
    @(-0.02446659572911941, 0.4848588817992852, 0.4654875719725148, 0.6822664962992943, 0.6891989838604905, 47),
    @(-0.1353368031892271, 0.4984219442182552, 0.4619043327933868, 0.6796354410957295, 0.6749645583256921, 38),
    @(-0.1049342653064903, 0.4821091839203969, 0.4458002174572263, 0.6676827221496946, 0.668480732802356, 37),
    @(-0.1747389245262045, 0.508949487189903, 0.4794270205088618, 0.6924066872213626, 0.6874003880098497, 20),
    @(-0.05712409056790357, 0.5421730072189278, 0.6097490892393119, 0.7808643219147049, 0.8105716633807137, 13),
    @(-0.5089197260409551, 0.5296605432532845, 0.5444046525602212, 0.7378378226685192, 0.597291140281082, 5)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($col = 2; $col -le 7; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 2]
    }
}

$wb.Save()
